$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 153
$ws.Cells.Item(6, 6).Value = 311
$ws.Cells.Item(7, 6).Value = 5481
$ws.Cells.Item(10, 6).Value = 3774
$ws.Cells.Item(11, 6).Value = 63
$ws.Cells.Item(14, 6).Value = 192
$ws.Cells.Item(17, 6).Value = 52
$ws.Cells.Item(20, 6).Value = 161
$ws.Cells.Item(23, 6).Value = 5151
$ws.Cells.Item(25, 6).Value = 2052
$ws.Cells.Item(27, 6).Value = 327
$ws.Cells.Item(28, 6).Value = 7631
$ws.Cells.Item(31, 6).Value = 2176
$ws.Cells.Item(32, 6).Value = 2139
$ws.Cells.Item(33, 6).Value = 1320
$ws.Cells.Item(34, 6).Value = 153
$ws.Cells.Item(35, 6).Value = 1165
$ws.Cells.Item(42, 6).Value = 1171
$ws.Cells.Item(44, 6).Value = 26
$ws.Cells.Item(45, 6).Value = 1301
$ws.Cells.Item(46, 6).Value = 2002
$ws.Cells.Item(47, 6).Value = 111
$ws.Cells.Item(48, 6).Value = 203
$ws.Cells.Item(49, 6).Value = 1205

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 6).Value = 145

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 6).Value = 720

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(3, 6).Value = 153
$ws.Cells.Item(7, 6).Value = 720
$ws.Cells.Item(8, 6).Value = 311
$ws.Cells.Item(9, 6).Value = 5481
$ws.Cells.Item(10, 6).Value = 3774
$ws.Cells.Item(11, 6).Value = 63
$ws.Cells.Item(14, 6).Value = 192
$ws.Cells.Item(16, 6).Value = 52
$ws.Cells.Item(19, 6).Value = 145
$ws.Cells.Item(20, 6).Value = 161
$ws.Cells.Item(24, 6).Value = 5151
$ws.Cells.Item(26, 6).Value = 2052
$ws.Cells.Item(28, 6).Value = 327
$ws.Cells.Item(29, 6).Value = 7631
$ws.Cells.Item(32, 6).Value = 2176
$ws.Cells.Item(33, 6).Value = 2139
$ws.Cells.Item(34, 6).Value = 1320
$ws.Cells.Item(35, 6).Value = 153
$ws.Cells.Item(36, 6).Value = 1165
$ws.Cells.Item(40, 6).Value = 1171
$ws.Cells.Item(42, 6).Value = 26
$ws.Cells.Item(43, 6).Value = 1301
$ws.Cells.Item(45, 6).Value = 2002
$ws.Cells.Item(46, 6).Value = 111
$ws.Cells.Item(48, 6).Value = 203
$ws.Cells.Item(49, 6).Value = 1205
